# feat(Marketing - Campañas): Agregadas variables a la importación de los contactos
#
# Adds six new "Variable N" columns (C:H) to the contact-import header row,
# right after the existing "NIT" / "Número de celular" columns, copies the
# header style onto them, widens the columns, and grows the header row
# height to fit the new wrapped headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header values -----------------------------------------------------
# A1/B1 ("NIT" / "Número de celular") are unchanged; insert the six new
# "Variable N" headers immediately after them.
$ws.Range("C1").Value = "Variable 1"
$ws.Range("D1").Value = "Variable 2"
$ws.Range("E1").Value = "Variable 3"
$ws.Range("F1").Value = "Variable 4"
$ws.Range("G1").Value = "Variable 5"
$ws.Range("H1").Value = "Variable 6"

# --- Match the header formatting used by A1/B1 on the new header cells -----
$ws.Range("A1").Copy()
$ws.Range("C1:H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Row height / column widths --------------------------------------------
$ws.Range("A1").EntireRow.RowHeight = 32.4

$ws.Range("A1").ColumnWidth = 27
$ws.Range("B1").ColumnWidth = 28.166666666666668
$ws.Range("C1").ColumnWidth = 23.666666666666668
$ws.Range("D1").ColumnWidth = 29
$ws.Range("E1").ColumnWidth = 30.833333333333332
$ws.Range("F1").ColumnWidth = 32.5
$ws.Range("G1").ColumnWidth = 32.833333333333336
$ws.Range("H1").ColumnWidth = 32.833333333333336

# --- Drop the stray leftover selection on A2:B2 -----------------------------
[void]$ws.Range("A1").Select()

Write-Host "Variables 1-6 added to the contact import header (C1:H1)."
